$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(325).Insert()

$ws.Cells.Item(325, 1).Value = 5
$ws.Cells.Item(325, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(325, 3).Value = "Maule"
$ws.Cells.Item(325, 4).Value = 44946
$ws.Cells.Item(325, 5).Value = 7
$ws.Cells.Item(325, 6).Value = 100112009
$ws.Cells.Item(325, 7).Value = "Acelga"
$ws.Cells.Item(325, 8).Value = "Sin especificar"
$ws.Cells.Item(325, 9).Value = "Primera"
$ws.Cells.Item(325, 10).Value = 200
$ws.Cells.Item(325, 11).Value = 3000
$ws.Cells.Item(325, 12).Value = 3000
$ws.Cells.Item(325, 13).Value = 3000
$ws.Cells.Item(325, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(325, 15).Value = "Región del Maule"
$ws.Cells.Item(325, 16).Value = 750
$ws.Cells.Item(325, 17).Value = 4
$ws.Cells.Item(325, 18).Value = "Hortaliza"
